$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "employee schema" fix drops two columns that were never used in
# practice: "dateOfJoin" (originally column I) and "todo" (originally
# column L). Deleting entire columns shifts everything after them left,
# which automatically re-numbers the remaining headers/data (basicSalary,
# montlyGoal, addOn, deduction) into their new positions (I, J, K, L).

# 1) Remove "dateOfJoin" (column I).
$ws.Columns("I:I").Delete() | Out-Null

# 2) Remove "todo". After the delete above it has shifted from column L
#    down to column K, so target K now.
$ws.Columns("K:K").Delete() | Out-Null

# The page was switched to portrait orientation.
$ws.PageSetup.Orientation = 1

# Restore the saved selection/active cell.
$ws.Range("J9").Select() | Out-Null
